$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 836
$ws.Cells.Item(28, 9).Value = 795
$ws.Cells.Item(28, 11).Value = 795
$ws.Cells.Item(28, 13).Value = -310

# ALC row 52
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(52, 8).Value = 499.66666
$ws.Cells.Item(52, 9).Value = 499
$ws.Cells.Item(52, 10).Value = 500
$ws.Cells.Item(52, 11).Value = 1497
$ws.Cells.Item(52, 12).Value = 1500
$ws.Cells.Item(52, 13).Value = -1337
$ws.Cells.Item(52, 14).Value = -1820

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 5600
$ws.Cells.Item(62, 9).Value = 4333.3335
$ws.Cells.Item(62, 10).Value = 7500
$ws.Cells.Item(62, 11).Value = 4333.3335
$ws.Cells.Item(62, 12).Value = 7500
$ws.Cells.Item(62, 13).Value = -3709.3335
$ws.Cells.Item(62, 14).Value = -8748

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 5600
$ws.Cells.Item(65, 9).Value = 4333.3335
$ws.Cells.Item(65, 10).Value = 7500
$ws.Cells.Item(65, 11).Value = 21666.6675
$ws.Cells.Item(65, 12).Value = 37500
$ws.Cells.Item(65, 13).Value = -18546.6675
$ws.Cells.Item(65, 14).Value = -43740

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 999.6667
$ws.Cells.Item(86, 10).Value = 999.6667
$ws.Cells.Item(86, 12).Value = 999.6667
$ws.Cells.Item(86, 14).Value = -3245.6667

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 999.6667
$ws.Cells.Item(89, 10).Value = 999.6667
$ws.Cells.Item(89, 12).Value = 4998.3335
$ws.Cells.Item(89, 14).Value = -16230.3335

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 3000
$ws.Cells.Item(113, 9).Value = 3000
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 3000
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 254
$ws.Cells.Item(113, 14).Value = ""

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1042
$ws.Cells.Item(2, 9).Value = 631.8570999999999
$ws.Cells.Item(2, 10).Value = 1999
$ws.Cells.Item(2, 11).Value = 631.8570999999999
$ws.Cells.Item(2, 12).Value = 1999
$ws.Cells.Item(2, 13).Value = -518.8570999999999
$ws.Cells.Item(2, 14).Value = -2225

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 401.42856
$ws.Cells.Item(97, 9).Value = 401.42856
$ws.Cells.Item(97, 11).Value = 401.42856
$ws.Cells.Item(97, 13).Value = 94.57144

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1042
$ws.Cells.Item(116, 9).Value = 631.8570999999999
$ws.Cells.Item(116, 10).Value = 1999
$ws.Cells.Item(116, 11).Value = 631.8570999999999
$ws.Cells.Item(116, 12).Value = 1999
$ws.Cells.Item(116, 13).Value = 1662.1429
$ws.Cells.Item(116, 14).Value = -6587

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3562.3333
$ws.Cells.Item(132, 9).Value = 3562.3333
$ws.Cells.Item(132, 11).Value = 10686.9999
$ws.Cells.Item(132, 13).Value = -8156.999899999999

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1042
$ws.Cells.Item(3, 9).Value = 631.8570999999999
$ws.Cells.Item(3, 10).Value = 1999
$ws.Cells.Item(3, 11).Value = 631.8570999999999
$ws.Cells.Item(3, 12).Value = 1999
$ws.Cells.Item(3, 13).Value = -517.8570999999999
$ws.Cells.Item(3, 14).Value = -2227

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 720.1429000000001
$ws.Cells.Item(80, 9).Value = 745.5
$ws.Cells.Item(80, 11).Value = 745.5
$ws.Cells.Item(80, 13).Value = 252.5

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value = 720.1429000000001
$ws.Cells.Item(83, 9).Value = 745.5
$ws.Cells.Item(83, 11).Value = 3727.5
$ws.Cells.Item(83, 13).Value = 1264.5

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 13).Value = ""

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 13).Value = ""

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 26.125
$ws.Cells.Item(7, 9).Value = 7.9
$ws.Cells.Item(7, 11).Value = 7.9
$ws.Cells.Item(7, 13).Value = 105.1

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 587.6
$ws.Cells.Item(22, 10).Value = 540
$ws.Cells.Item(22, 12).Value = 540
$ws.Cells.Item(22, 14).Value = -1240

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 900
$ws.Cells.Item(58, 9).Value = 900
$ws.Cells.Item(58, 11).Value = 900
$ws.Cells.Item(58, 13).Value = -697

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 900
$ws.Cells.Item(136, 9).Value = 900
$ws.Cells.Item(136, 11).Value = 2700
$ws.Cells.Item(136, 13).Value = -150

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1100
$ws.Cells.Item(4, 9).Value = 1500
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 4500
$ws.Cells.Item(4, 12).Value = 900
$ws.Cells.Item(4, 13).Value = -4388
$ws.Cells.Item(4, 14).Value = -1124

# CUL row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 1000
$ws.Cells.Item(22, 9).Value = 1000
$ws.Cells.Item(22, 11).Value = 3000
$ws.Cells.Item(22, 13).Value = -2831

# CUL row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(27, 8).Value = 1000
$ws.Cells.Item(27, 9).Value = 1000
$ws.Cells.Item(27, 11).Value = 3000
$ws.Cells.Item(27, 13).Value = -2898

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1750
$ws.Cells.Item(34, 9).Value = 500
$ws.Cells.Item(34, 11).Value = 1500
$ws.Cells.Item(34, 13).Value = -1416

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 1556.3572
$ws.Cells.Item(39, 10).Value = 2429.875
$ws.Cells.Item(39, 12).Value = 7289.625
$ws.Cells.Item(39, 14).Value = -7877.625

# CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 25
$ws.Cells.Item(46, 9).Value = 25
$ws.Cells.Item(46, 11).Value = 75
$ws.Cells.Item(46, 13).Value = 16

# CUL row 49
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(49, 8).Value = 1273
$ws.Cells.Item(49, 10).Value = 841.5
$ws.Cells.Item(49, 12).Value = 2524.5
$ws.Cells.Item(49, 14).Value = -2836.5

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1348
$ws.Cells.Item(97, 9).Value = 1121.25
$ws.Cells.Item(97, 10).Value = 3162
$ws.Cells.Item(97, 11).Value = 1121.25
$ws.Cells.Item(97, 12).Value = 3162
$ws.Cells.Item(97, 13).Value = -625.25
$ws.Cells.Item(97, 14).Value = -4154

# GSM row 104
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(104, 8).Value = 132200
$ws.Cells.Item(104, 10).Value = 132200
$ws.Cells.Item(104, 12).Value = 132200
$ws.Cells.Item(104, 14).Value = -139188

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2749.25
$ws.Cells.Item(132, 9).Value = 2166.1667
$ws.Cells.Item(132, 10).Value = 4498.5
$ws.Cells.Item(132, 11).Value = 6498.500100000001
$ws.Cells.Item(132, 12).Value = 13495.5
$ws.Cells.Item(132, 13).Value = -3968.500100000001
$ws.Cells.Item(132, 14).Value = -18555.5

# LTW row 74
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(74, 8).Value = 52000
$ws.Cells.Item(74, 10).Value = 52000
$ws.Cells.Item(74, 12).Value = 52000
$ws.Cells.Item(74, 14).Value = -53996

# LTW row 77
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(77, 8).Value = 52000
$ws.Cells.Item(77, 10).Value = 52000
$ws.Cells.Item(77, 12).Value = 156000
$ws.Cells.Item(77, 14).Value = -165984

# WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 75000
$ws.Cells.Item(46, 10).Value = 75000
$ws.Cells.Item(46, 12).Value = 75000
$ws.Cells.Item(46, 14).Value = -75462

# WVR row 61
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(61, 8).Value = 25664.666
$ws.Cells.Item(61, 9).Value = 25664.666
$ws.Cells.Item(61, 11).Value = 25664.666
$ws.Cells.Item(61, 13).Value = -25372.666

# WVR row 95
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).Value = ""

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1114.6666
$ws.Cells.Item(107, 10).Value = 1486
$ws.Cells.Item(107, 12).Value = 4458
$ws.Cells.Item(107, 14).Value = -8298

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 17832.5
$ws.Cells.Item(122, 9).Value = 1665.6666
$ws.Cells.Item(122, 11).Value = 4996.9998
$ws.Cells.Item(122, 13).Value = -2546.9998

# WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(134, 8).Value = 75000
$ws.Cells.Item(134, 10).Value = 75000
$ws.Cells.Item(134, 12).Value = 225000
$ws.Cells.Item(134, 14).Value = -230070

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(135, 8).Value = 41666.332
$ws.Cells.Item(135, 9).Value = 34999
$ws.Cells.Item(135, 10).Value = 45000
$ws.Cells.Item(135, 11).Value = 34999
$ws.Cells.Item(135, 12).Value = 45000
$ws.Cells.Item(135, 13).Value = -29929
$ws.Cells.Item(135, 14).Value = -55140
